# Expense-tracker workbook update:
#  - fix mojibake in the "AtB ... monedskort" transportation entry
#  - remove the stray test row ("asdfg")
#  - add the new transactions that were entered using the new
#    "press Enter to submit" flow (pizza, buss, sdaf)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("March")

# --- fix the mis-encoded Norwegian text in row 6 ---
$ws.Range("B6").Value = "AtB månedskort"

# --- drop the old "asdfg" test entry (row 7); HYRE shifts up from row 8 ---
$ws.Rows(7).Delete()

# --- helper range used to stamp "plain" (unstyled) formatting onto new cells ---
$plainFormat = $ws.Range("C2:D2")

function Add-Transaction {
    param(
        [int]$row,
        [string]$category,
        [string]$name,
        [string]$date,
        [string]$price,
        [string]$account
    )

    $ws.Cells.Item($row, 1).Value = $category
    $ws.Cells.Item($row, 2).Value = $name

    # Date/Price columns must stay plain text (matching the rest of the
    # sheet) instead of being auto-converted to a date/number, so mark
    # them as text before typing the values.
    $dateCell = $ws.Cells.Item($row, 3)
    $priceCell = $ws.Cells.Item($row, 4)
    $dateCell.NumberFormat = "@"
    $priceCell.NumberFormat = "@"
    $dateCell.Value = $date
    $priceCell.Value = $price

    $ws.Cells.Item($row, 5).Value = $account

    # Re-apply the same (default) formatting used elsewhere on the sheet
    # so the new cells don't carry a leftover text-format style.
    $plainFormat.Copy()
    $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 4)).PasteSpecial(-4122)
    $excel.CutCopyMode = 0
}

Add-Transaction 8 "Food" "pizza" "2023-03-09" "1000.0" "Checkings"
Add-Transaction 9 "Transportation" "buss" "2023-03-18" "2000.0" "Savings"
Add-Transaction 10 "Food" "sdaf" "2023-03-18" "1234.0" "Checkings"
